# Delete two slides from the presentation:
#  - the original 1st slide (SlideID 256)
#  - the original 3rd slide (SlideID 257)
# After removing the 1st slide, the former 3rd slide shifts to position 2,
# so we delete index 1 then index 2.

$p = $ppt.ActivePresentation

$p.Slides.Item(1).Delete()
$p.Slides.Item(2).Delete()
